# Update weekly Fruta/Hortaliza (Frambuesa - Vega Monumental Concepción) price data:
# dates, calidad (Primera/Segunda), volumen and min/max/avg/$-per-kg prices per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44574
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7500
$ws.Range("S2").Value = 3750
$ws.Range("D3").Value = 44574
$ws.Range("L3").Value = "Segunda"
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000
$ws.Range("D4").Value = 44195
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3000
$ws.Range("O4").Value = 3500
$ws.Range("P4").Value = 3250
$ws.Range("S4").Value = 1625
$ws.Range("D5").Value = 44195
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 2500
$ws.Range("O5").Value = 2500
$ws.Range("P5").Value = 2500
$ws.Range("S5").Value = 1250
$ws.Range("D6").Value = 44602
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6500
$ws.Range("S6").Value = 3250
$ws.Range("D7").Value = 44602
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("S7").Value = 2500
$ws.Range("D8").Value = 44609
$ws.Range("L8").Value = "Primera"
$ws.Range("N8").Value = 6500
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 6750
$ws.Range("S8").Value = 3375
$ws.Range("D9").Value = 44609
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("S9").Value = 3000
$ws.Range("D10").Value = 44532
$ws.Range("L10").Value = "Primera"
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 5000
$ws.Range("D11").Value = 44532
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 4000
$ws.Range("D12").Value = 44617
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6500
$ws.Range("S12").Value = 3250
$ws.Range("D13").Value = 44559
$ws.Range("D14").Value = 44559
$ws.Range("D15").Value = 44216
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 3500
$ws.Range("O15").Value = 4000
$ws.Range("P15").Value = 3750
$ws.Range("S15").Value = 1875
$ws.Range("D16").Value = 44216
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3000
$ws.Range("P16").Value = 3000
$ws.Range("S16").Value = 1500
